# Auto-generated edit script to apply market-data refresh values
# as described by the commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3068
$ws.Range("I98").Value = 2826.6667
$ws.Range("J98").Value = 3309.3333
$ws.Range("K98").Value = 2826.6667
$ws.Range("L98").Value = 3309.3333
$ws.Range("M98").Value = -1328.6667
$ws.Range("N98").Value = -6305.3333

$ws.Range("H122").Value = 3068
$ws.Range("I122").Value = 2826.6667
$ws.Range("J122").Value = 3309.3333
$ws.Range("K122").Value = 8480.000100000001
$ws.Range("L122").Value = 9927.999899999999
$ws.Range("M122").Value = -6030.000100000001
$ws.Range("N122").Value = -14827.9999

$ws.Range("H137").Value = 1686.8889
$ws.Range("I137").Value = 1193
$ws.Range("J137").Value = 2082
$ws.Range("K137").Value = 3579
$ws.Range("L137").Value = 6246
$ws.Range("M137").Value = -1029
$ws.Range("N137").Value = -11346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1983.4166
$ws.Range("I61").Value = 1150.3334
$ws.Range("J61").Value = 2816.5
$ws.Range("K61").Value = 1150.3334
$ws.Range("L61").Value = 2816.5
$ws.Range("M61").Value = -938.3334
$ws.Range("N61").Value = -3240.5

$ws.Range("H74").Value = 1038.875
$ws.Range("I74").Value = 1038.875
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1038.875
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -164.875
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1038.875
$ws.Range("I77").Value = 1038.875
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5194.375
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -826.375
$ws.Range("N77").ClearContents()

$ws.Range("H132").Value = 2248.2
$ws.Range("I132").Value = 1347.08
$ws.Range("J132").Value = 3750.0667
$ws.Range("K132").Value = 4041.24
$ws.Range("L132").Value = 11250.2001
$ws.Range("M132").Value = -1511.24
$ws.Range("N132").Value = -16310.2001

$ws.Range("H136").Value = 1983.4166
$ws.Range("I136").Value = 1150.3334
$ws.Range("J136").Value = 2816.5
$ws.Range("K136").Value = 3451.0002
$ws.Range("L136").Value = 8449.5
$ws.Range("M136").Value = -901.0001999999999
$ws.Range("N136").Value = -13549.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18536.092
$ws.Range("I82").Value = 8909.5
$ws.Range("J82").Value = 30088
$ws.Range("K82").Value = 8909.5
$ws.Range("L82").Value = 30088
$ws.Range("M82").Value = -8526.5
$ws.Range("N82").Value = -30854

$ws.Range("H85").Value = 18536.092
$ws.Range("I85").Value = 8909.5
$ws.Range("J85").Value = 30088
$ws.Range("K85").Value = 8909.5
$ws.Range("L85").Value = 30088
$ws.Range("M85").Value = -7583.5
$ws.Range("N85").Value = -32740

$ws.Range("H86").Value = 1655.2222
$ws.Range("I86").Value = 1419.6
$ws.Range("J86").Value = 1949.75
$ws.Range("K86").Value = 1419.6
$ws.Range("L86").Value = 1949.75
$ws.Range("M86").Value = -296.5999999999999
$ws.Range("N86").Value = -4195.75

$ws.Range("H89").Value = 1655.2222
$ws.Range("I89").Value = 1419.6
$ws.Range("J89").Value = 1949.75
$ws.Range("K89").Value = 7098
$ws.Range("L89").Value = 9748.75
$ws.Range("M89").Value = -1482
$ws.Range("N89").Value = -20980.75

$ws.Range("H99").Value = 1407.7826
$ws.Range("I99").Value = 1348.8948
$ws.Range("J99").Value = 1687.5
$ws.Range("K99").Value = 1348.8948
$ws.Range("L99").Value = 1687.5
$ws.Range("M99").Value = 149.1052
$ws.Range("N99").Value = -4683.5

$ws.Range("H134").Value = 2159.739
$ws.Range("I134").Value = 2016.1177
$ws.Range("J134").Value = 2566.6667
$ws.Range("K134").Value = 6048.3531
$ws.Range("L134").Value = 7700.000100000001
$ws.Range("M134").Value = -3513.3531
$ws.Range("N134").Value = -12770.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2771.9268
$ws.Range("I31").Value = 2305.04
$ws.Range("J31").Value = 3501.4375
$ws.Range("K31").Value = 2305.04
$ws.Range("L31").Value = 3501.4375
$ws.Range("M31").Value = -2010.04
$ws.Range("N31").Value = -4091.4375

$ws.Range("H34").Value = 2771.9268
$ws.Range("I34").Value = 2305.04
$ws.Range("J34").Value = 3501.4375
$ws.Range("K34").Value = 2305.04
$ws.Range("L34").Value = 3501.4375
$ws.Range("M34").Value = -2103.04
$ws.Range("N34").Value = -3905.4375

$ws.Range("H58").Value = 1198.591
$ws.Range("I58").Value = 825.4545000000001
$ws.Range("J58").Value = 1571.7273
$ws.Range("K58").Value = 825.4545000000001
$ws.Range("L58").Value = 1571.7273
$ws.Range("M58").Value = -622.4545000000001
$ws.Range("N58").Value = -1977.7273

$ws.Range("H132").Value = 6059
$ws.Range("I132").Value = 6328.3335
$ws.Range("J132").Value = 5116.3335
$ws.Range("K132").Value = 18985.0005
$ws.Range("L132").Value = 15349.0005
$ws.Range("M132").Value = -16455.0005
$ws.Range("N132").Value = -20409.0005

$ws.Range("H134").Value = 1889.3158
$ws.Range("I134").Value = 1354.9286
$ws.Range("J134").Value = 3385.6
$ws.Range("K134").Value = 4064.7858
$ws.Range("L134").Value = 10156.8
$ws.Range("M134").Value = -1529.7858
$ws.Range("N134").Value = -15226.8

$ws.Range("H136").Value = 1198.591
$ws.Range("I136").Value = 825.4545000000001
$ws.Range("J136").Value = 1571.7273
$ws.Range("K136").Value = 2476.3635
$ws.Range("L136").Value = 4715.1819
$ws.Range("M136").Value = 73.63649999999961
$ws.Range("N136").Value = -9815.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 31990390
$ws.Range("I121").Value = 350
$ws.Range("J121").Value = 36402810
$ws.Range("K121").Value = 1050
$ws.Range("L121").Value = 109208430
$ws.Range("M121").Value = 260
$ws.Range("N121").Value = -109211050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4301.077
$ws.Range("I102").Value = 4000
$ws.Range("J102").Value = 4434.8887
$ws.Range("K102").Value = 4000
$ws.Range("L102").Value = 4434.8887
$ws.Range("M102").Value = -2378
$ws.Range("N102").Value = -7678.8887

$ws.Range("H132").Value = 4822
$ws.Range("I132").Value = 6371
$ws.Range("J132").Value = 2438.923
$ws.Range("K132").Value = 19113
$ws.Range("L132").Value = 7316.768999999999
$ws.Range("M132").Value = -16583
$ws.Range("N132").Value = -12376.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 281.69232
$ws.Range("I55").Value = 251.66667
$ws.Range("J55").Value = 307.42856
$ws.Range("K55").Value = 251.66667
$ws.Range("L55").Value = 307.42856
$ws.Range("M55").Value = -78.66667000000001
$ws.Range("N55").Value = -653.4285600000001

$ws.Range("H132").Value = 2950.3845
$ws.Range("I132").Value = 1913.5333
$ws.Range("J132").Value = 4364.273
$ws.Range("K132").Value = 5740.5999
$ws.Range("L132").Value = 13092.819
$ws.Range("M132").Value = -3210.5999
$ws.Range("N132").Value = -18152.819

$ws.Range("H136").Value = 2669.1765
$ws.Range("I136").Value = 1483.3334
$ws.Range("J136").Value = 3316
$ws.Range("K136").Value = 4450.0002
$ws.Range("L136").Value = 9948
$ws.Range("M136").Value = -1900.0002
$ws.Range("N136").Value = -15048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 27623.46
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 27623.46
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 27623.46
$ws.Range("N70").Value = -28253.46

$ws.Range("H73").Value = 27623.46
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 27623.46
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 27623.46
$ws.Range("N73").Value = -29807.46

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H122").Value = 5342.7144
$ws.Range("I122").Value = 5079.8
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 15239.4
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -12789.4
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 1674.3235
$ws.Range("I132").Value = 876
$ws.Range("J132").Value = 2168.524
$ws.Range("K132").Value = 2628
$ws.Range("L132").Value = 6505.572
$ws.Range("M132").Value = -98
$ws.Range("N132").Value = -11565.572

$ws.Range("H136").Value = 1331.85
$ws.Range("I136").Value = 1036.8334
$ws.Range("J136").Value = 1774.375
$ws.Range("K136").Value = 3110.5002
$ws.Range("L136").Value = 5323.125
$ws.Range("M136").Value = -560.5001999999999
$ws.Range("N136").Value = -10423.125
